# 因卓教育项目计划 - 时间进度更新 (yin zhuo jiao yu time schedule)
#
# Updates the source-data table (start/finish dates + estimated workload),
# resizes / restyles the Gantt-style bar chart on the "图标" sheet, and
# brings the two sheet views back to where the author left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 源数据 sheet - task schedule table
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

# Row 3 - 原型设计
$ws1.Range("B3").Value = 42853
$ws1.Range("C3").Formula = "=20.035+9.05"
$ws1.Range("E3").Value = 42902

# Row 4 - UI设计
$ws1.Range("B4").Value = 42906
$ws1.Range("C4").Formula = "=22.63 + 18.88"
$ws1.Range("E4").Value = 43000

# Row 5 - iOS APP开发
$ws1.Range("B5").Value = 43003
$ws1.Range("C5").Value = 62.5
$ws1.Range("E5").Value = 43065

# Row 6 - Android APP 开发
$ws1.Range("B6").Value = 43003
$ws1.Range("C6").Value = 62.5
$ws1.Range("E6").Value = 43065

# Row 7 - PC 前端 开发
$ws1.Range("B7").Value = 43003
$ws1.Range("C7").Value = 58.5
$ws1.Range("E7").Value = 43063

# Row 8 - 后台
$ws1.Range("B8").Value = 42906
$ws1.Range("C8").Formula = "=80.38"
$ws1.Range("E8").Value = 43031

# E8 previously carried its own one-off number format (General, flagged as
# applied); bring it back in line with the rest of column E (m/d/yyyy),
# matching E3's style exactly so it shares the same style record.
$ws1.Range("E3").Copy()
$ws1.Range("E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Calculate()

# restore the author's cell selection on this sheet
$ws1.Range("B3").Select()

# ---------------------------------------------------------------------
# 2. 图标 sheet - chart: resize, recolor, and fix the view
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$co = $ws2.ChartObjects().Item(1)
$chart = $co.Chart

# Drag-resize the chart frame (anchored from A1) out to roughly col T / row 45.
$co.Width = 1156.3839173228346
$co.Height = 630.7174015748031

# Recolor the "耗费天数" series from accent2 to a darkened accent1.
$series2 = $chart.SeriesCollection().Item(2)
$series2.Format.Fill.Solid()
$series2.Format.Fill.ForeColor.ObjectThemeColor = 5
$series2.Format.Fill.ForeColor.TintAndShade = 0
$series2.Format.Fill.ForeColor.Brightness = -0.5

# Lighten the axis label text color (drop the 65%/35% lumMod/lumOff tx1 tint
# back to plain tx1) on both axes.
$catAxis = $chart.Axes(1)
$valAxis = $chart.Axes(2)
$catAxis.TickLabels.Font.ThemeColor = 1
$valAxis.TickLabels.Font.ThemeColor = 1

$ws2.Range("W22").Select()
$excel.ActiveWindow.Zoom = 70

$wb.Save()
